# Feature tracker update:
#  - Add "Never" as the Completed Version for the two dice-system rows
#    (Add Genesys dice / Add Fantasy Flight Star Wars dice).
#  - Add three new feature requests from Weston Fiala at the bottom of
#    the tracker (rows 27-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows.
$ws.Range("A27").Value = "Quick category change"
$ws.Range("B27").Value = "I want to be able to change my saved roll to have a different existing category easily"
$ws.Range("D27").Value = "Weston Fiala"

$ws.Range("A28").Value = "Compact view - Saved"
$ws.Range("A29").Value = "Compact view - Custom"
$ws.Range("D28").Value = "Weston Fiala"
$ws.Range("D29").Value = "Weston Fiala"

$ws.Range("B29").Value = "Custom roll items take up too much space each, I want more to fit on one screen."
$ws.Range("B28").Value = "Saved roll items take up too much space each, I want more to fit on one screen."

# Completed Version ("Never") for the two existing dice rows.
$ws.Range("C17").Value = "Never"
$ws.Range("C18").Value = "Never"

# Column A needs to widen to fit the new, longer labels.
$ws.Columns("A").AutoFit()

# Leave the cursor where the author's session ended up.
$ws.Range("A32").Select()
